# update scripts wuth new tpm
# Refresh the NATMI ligand/receptor expression + specificity + edge-weight
# figures on Sheet1 for the Adam17-Notch1 LR pair using the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.771609
$ws.Range("H2").Value = 26.314827
$ws.Range("I2").Value = 0.2200338127677125
$ws.Range("J2").Value = 0.2200338127677125
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 338.1690056234233
$ws.Range("R2").Value = 3043.52105061081
$ws.Range("S2").Value = 0.1266988078740261
$ws.Range("T2").Value = 0.1266988078740261
$ws.Range("G3").Value = 8.771609
$ws.Range("H3").Value = 26.314827
$ws.Range("I3").Value = 0.2200338127677125
$ws.Range("J3").Value = 0.2200338127677125
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 47.109469937684
$ws.Range("R3").Value = 423.9852294391561
$ws.Range("S3").Value = 0.01765009087594635
$ws.Range("T3").Value = 0.01765009087594635
$ws.Range("G4").Value = 8.771609
$ws.Range("H4").Value = 26.314827
$ws.Range("I4").Value = 0.2200338127677125
$ws.Range("J4").Value = 0.2200338127677125
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 202.0089418640883
$ws.Range("R4").Value = 1818.080476776795
$ws.Range("S4").Value = 0.07568491401774001
$ws.Range("T4").Value = 0.07568491401774004
$ws.Range("I5").Value = 0.583164828467109
$ws.Range("J5").Value = 0.583164828467109
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 896.2634773114045
$ws.Range("R5").Value = 8066.371295802641
$ws.Range("S5").Value = 0.3357951563510134
$ws.Range("T5").Value = 0.3357951563510134
$ws.Range("I6").Value = 0.583164828467109
$ws.Range("J6").Value = 0.583164828467109
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("S6").Value = 0.0467787749920339
$ws.Range("T6").Value = 0.04677877499203391
$ws.Range("I7").Value = 0.583164828467109
$ws.Range("J7").Value = 0.583164828467109
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 535.3927582728312
$ws.Range("R7").Value = 4818.534824455481
$ws.Range("S7").Value = 0.2005908971240616
$ws.Range("T7").Value = 0.2005908971240617
$ws.Range("G8").Value = 7.845451333333334
$ws.Range("H8").Value = 23.536354
$ws.Range("I8").Value = 0.1968013587651783
$ws.Range("J8").Value = 0.1968013587651783
$ws.Range("M8").Value = 38.55267666666666
$ws.Range("N8").Value = 115.65803
$ws.Range("O8").Value = 0.5758151725879548
$ws.Range("P8").Value = 0.5758151725879548
$ws.Range("Q8").Value = 302.4631485580689
$ws.Range("R8").Value = 2722.16833702262
$ws.Range("S8").Value = 0.1133212083629152
$ws.Range("T8").Value = 0.1133212083629152
$ws.Range("G9").Value = 7.845451333333334
$ws.Range("H9").Value = 23.536354
$ws.Range("I9").Value = 0.1968013587651783
$ws.Range("J9").Value = 0.1968013587651783
$ws.Range("O9").Value = 0.08021535714867321
$ws.Range("P9").Value = 0.08021535714867323
$ws.Range("Q9").Value = 42.13537718510134
$ws.Range("R9").Value = 379.2183946659121
$ws.Range("S9").Value = 0.01578649128069295
$ws.Range("T9").Value = 0.01578649128069295
$ws.Range("G10").Value = 7.845451333333334
$ws.Range("H10").Value = 23.536354
$ws.Range("I10").Value = 0.1968013587651783
$ws.Range("J10").Value = 0.1968013587651783
$ws.Range("M10").Value = 23.02986166666667
$ws.Range("N10").Value = 69.089585
$ws.Range("O10").Value = 0.3439694702633719
$ws.Range("P10").Value = 0.3439694702633719
$ws.Range("Q10").Value = 180.6796589192322
$ws.Range("R10").Value = 1626.11693027309
$ws.Range("S10").Value = 0.06769365912157019
$ws.Range("T10").Value = 0.0676936591215702
